# Apply updated dSF (column F) values for specific rows.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @{
    5  = -7
    10 = -5
    20 = -8
    22 = -9
    23 = -5
    24 = -4
    25 = 9
    28 = 0
    30 = -7
    31 = -1
    33 = -8
    35 = -4
    40 = -5
    42 = 2
}

foreach ($row in $updates.Keys) {
    $ws.Range("F$row").Value = $updates[$row]
}
